$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "Salade de pomme de terre"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("G3").Value = "Crépinette de porc BBC"
$ws.Range("H3").Value = "Boulettes pois chiches"
$ws.Range("I3").Value = "Lentilles"
$ws.Range("J3").Value = "Carottes braisées"
$ws.Range("K3").Value = ""

$ws.Range("E5").Select()
